$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# NOTE: all offsets below are 1-based character positions measured against
# the ORIGINAL (unedited) TextRange.Text. Edits are applied from the
# highest original offset down to the lowest so earlier offsets remain
# valid while later text is being shortened/lengthened.

# --- "Sao Carlos, {{DATA}}" paragraph: fold the ", " run into "Sao Carlos" ---
$cDataSep = $tr.Characters(225, 2)
$cDataSep.Text = ""

$cSaoCarlos = $tr.Characters(215, 10)
$cSaoCarlos.Text = "São Carlos, "

# --- "portador ... RG ... CPF" paragraph: drop the RG mention ---
# Replace "e CPF nº " with "do CPF nº " first, then re-split that new text
# into three runs ("do " / "CPF " / "nº ") by re-writing each sub-piece
# (right to left so positions stay valid).
$cEcpf = $tr.Characters(53, 9)
$cEcpf.Text = "do CPF nº "

$cNo = $tr.Characters(60, 3)
$cNo.Text = "nº "

$cCpf = $tr.Characters(56, 4)
$cCpf.Text = "CPF "

$cDo = $tr.Characters(53, 3)
$cDo.Text = "do "

# Remove the "{{RG}} " placeholder run entirely.
$cRg = $tr.Characters(46, 7)
$cRg.Text = ""

# "portador do RG nº " -> "portador "
$cPortador = $tr.Characters(28, 18)
$cPortador.Text = "portador "
